{"js": "// The original document only contained the placeholder text \"JD\" plus a\n// leftover \"_GoBack\" bookmark from the last cursor position. Give the\n// sample document some real content: replace the placeholder text with a\n// descriptive sentence and drop the stale bookmark.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Replace the paragraph's text in place (keeps the paragraph itself,\n// just swaps the run content).\nfirstParagraph\n  .getRange()\n  .insertText(\n    \"This is a sample document that can be converted using office-converter\",\n    Word.InsertLocation.replace\n  );\n\n// Remove the stale \"_GoBack\" bookmark; it's a no-op if it isn't present.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Replace the sample placeholder text \"JD\" with the full sample sentence,\n# and remove the leftover \"_GoBack\" bookmark (Word housekeeping artifact\n# that the re-saved document no longer contains).\n\n$d = $word.ActiveDocument\n\n$firstParagraph = $d.Paragraphs(1)\n$firstParagraph.Range.Text = \"This is a sample document that can be converted using office-converter\"\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$d.Save()\n"}
